$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5496483445167542
$ws.Range("B1").Value = 1.077258825302124
$ws.Range("C1").Value = 5.185221672058105
$ws.Range("D1").Value = 3.982270240783691
$ws.Range("E1").Value = 1.056949496269226
